$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.494.30'
$ws.Range('E2').Value = '  +1.94%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.863.94'
$ws.Range('E3').Value = '  +0.90%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.010'
$ws.Range('E4').Value = '  -0.38%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.42'
$ws.Range('E5').Value = '  +0.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.009'
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4778'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3815'
$ws.Range('E8').Value = '  +3.77%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07335'
$ws.Range('E9').Value = '  +1.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9362'
$ws.Range('E10').Value = '  +0.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.80'
$ws.Range('E11').Value = '  +5.61%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07809'
$ws.Range('E12').Value = '  +0.93%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.865.15'
$ws.Range('E13').Value = '  +5.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.446'
$ws.Range('E14').Value = '  +1.85%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.555'
$ws.Range('E15').Value = '  +1.93%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '90.53'
$ws.Range('E16').Value = '  +2.00%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.012'
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008785'
$ws.Range('E18').Value = '  +1.65%  '
$ws.Range('E19').Value = '  -0.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.534.75'
$ws.Range('E20').Value = '  +2.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.67'
$ws.Range('E21').Value = '  +1.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.122'
$ws.Range('E22').Value = '  +1.23%  '
$ws.Range('E23').Value = '  +0.59%  '
$ws.Range('E24').Value = '  +0.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.72'
$ws.Range('E25').Value = '  +1.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.49'
$ws.Range('E26').Value = '  +1.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.019'
$ws.Range('E27').Value = '  +1.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '115.34'
$ws.Range('E28').Value = '  +1.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.939'
$ws.Range('E29').Value = '  -0.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08892'
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.320'
$ws.Range('E31').Value = '  -0.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.218'
$ws.Range('E32').Value = '  +4.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7589'
$ws.Range('E33').Value = '  +2.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.605'
$ws.Range('E34').Value = '  +2.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.726'
$ws.Range('E35').Value = '  -0.72%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02054'
$ws.Range('E36').Value = '  +4.72%  '
$ws.Range('E37').Value = '  +0.37%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5573'
$ws.Range('E38').Value = '  +7.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05279'
$ws.Range('E39').Value = '  +0.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.992'
$ws.Range('E40').Value = '  +0.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.062'
$ws.Range('E41').Value = '  +1.09%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.669'
$ws.Range('E42').Value = '  +5.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1526'
$ws.Range('E43').Value = '  +0.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4909'
$ws.Range('E44').Value = '  +3.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.66'
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.010'
$ws.Range('E46').Value = '  -0.28%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.657'
$ws.Range('E47').Value = '  +3.02%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '103.04'
$ws.Range('E48').Value = '  +1.36%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '67.45'
$ws.Range('E49').Value = '  +3.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06089'
$ws.Range('E50').Value = '  +0.34%  '
$ws.Range('E51').Value = '  +3.10%  '
